# Ajout EDT S6 25-26.
# "Securite" and "Legislation" course codes gain a second linked group
# (KUPT9BB1 / KUPT9BA2), and the room column (F) is filled in for the
# lecture/TD sessions that previously had no room assigned.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Sécurité - C (KRTA9AA3/KUPT9BB1)"
$ws.Range("B3").Value = "KRTA9AA3/KUPT9BB1"
$ws.Range("F3").Value = "U3-110"
$ws.Range("A6").Value = "Sécurité - C (KRTA9AA3/KUPT9BB1)"
$ws.Range("B6").Value = "KRTA9AA3/KUPT9BB1"
$ws.Range("F6").Value = "U3-110"
$ws.Range("A7").Value = "Sécurité - C (KRTA9AA3/KUPT9BB1)"
$ws.Range("B7").Value = "KRTA9AA3/KUPT9BB1"
$ws.Range("F7").Value = "U3-110"
$ws.Range("A10").Value = "Sécurité - TP (KRTA9AA3/KUPT9BB1) (FC)"
$ws.Range("B10").Value = "KRTA9AA3/KUPT9BB1"
$ws.Range("A11").Value = "Sécurité - TP (KRTA9AA3/KUPT9BB1) (FC)"
$ws.Range("B11").Value = "KRTA9AA3/KUPT9BB1"
$ws.Range("A14").Value = "Législation - C (KRTA9AD1/KUPT9BA2)"
$ws.Range("B14").Value = "KRTA9AD1/KUPT9BA2"
$ws.Range("F14").Value = "U3-110"
$ws.Range("A17").Value = "Législation - C (KRTA9AD1/KUPT9BA2)"
$ws.Range("B17").Value = "KRTA9AD1/KUPT9BA2"
$ws.Range("F17").Value = "U3-110"
$ws.Range("A18").Value = "Législation - C (KRTA9AD1/KUPT9BA2)"
$ws.Range("B18").Value = "KRTA9AD1/KUPT9BA2"
$ws.Range("F18").Value = "U3-110"
$ws.Range("A21").Value = "Législation - C (KRTA9AD1/KUPT9BA2)"
$ws.Range("B21").Value = "KRTA9AD1/KUPT9BA2"
$ws.Range("F21").Value = "U3-4"
$ws.Range("A24").Value = "Législation - TD (KRTA9AD1/KUPT9BA2)"
$ws.Range("B24").Value = "KRTA9AD1/KUPT9BA2"
$ws.Range("F24").Value = "U3-Amphi"
$ws.Range("A26").Value = "Législation - TD (KRTA9AD1/KUPT9BA2)"
$ws.Range("B26").Value = "KRTA9AD1/KUPT9BA2"
$ws.Range("F26").Value = "U3-Amphi"
$ws.Range("A27").Value = "Législation - TD (KRTA9AD1/KUPT9BA2)"
$ws.Range("B27").Value = "KRTA9AD1/KUPT9BA2"
$ws.Range("F27").Value = "U3-Amphi"
$ws.Range("A29").Value = "Législation - TD (KRTA9AD1/KUPT9BA2)"
$ws.Range("B29").Value = "KRTA9AD1/KUPT9BA2"
$ws.Range("F29").Value = "U3-Amphi"
$ws.Range("F32").Value = "U3-Amphi"
$ws.Range("F39").Value = "U3-Amphi"
$ws.Range("F40").Value = "U3-Amphi"
$ws.Range("F41").Value = "U3-Amphi"
$ws.Range("A44").Value = "Législation - C (KRTA9AD1/KUPT9BA2)"
$ws.Range("B44").Value = "KRTA9AD1/KUPT9BA2"
$ws.Range("F44").Value = "U3-Amphi"
$ws.Range("F46").Value = "U3-Amphi"
